# Regenerate merged AHB files
# 1) Rename the header row: "<name>_old" -> "<name>_FV2310", "<name>_new" -> "<name>_FV2404"
# 2) Turn the data range into an Excel Table ("Table1")
# 3) Freeze the header row (pane split after row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 21   # A .. U
$lastRow = 84

for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value()
    if ($val -like "*_old") {
        $cell.Value = ($val -replace "_old$", "_FV2310")
    } elseif ($val -like "*_new") {
        $cell.Value = ($val -replace "_new$", "_FV2404")
    }
}

# Convert the range to a proper Excel Table (ListObject), keeping the header
# names already written above.
$tableRange = $ws.Range("A1:U84")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# Freeze panes above row 2 (i.e. freeze the header row).
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
